# Refresh the cryptos table with the latest scraped price/volume snapshot.
# (GitHub Actions nightly data refresh.)
#
# Price-column values are plain scraped text (e.g. "0.9993", "28.840.09") and
# must stay text, not become numbers. Cells whose text parses as a plain number
# are written with a leading apostrophe, same as typing '0.9993 into Excel,
# so the stored cell stays a text value instead of being coerced to a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.840.09'
$ws.Range('E2').Value = '  +7.77%  '

$ws.Range('D3').Value = '1.811.26'
$ws.Range('E3').Value = '  +4.97%  '

$ws.Range('D4').Value = '''0.9993'
$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').Value = '''249.25'
$ws.Range('E5').Value = '  +3.56%  '

$ws.Range('D6').Value = '''0.9996'
$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('E7').Value = '  +1.80%  '

$ws.Range('D8').Value = '''0.2785'
$ws.Range('E8').Value = '  +7.54%  '

$ws.Range('D9').Value = '''0.06447'
$ws.Range('E9').Value = '  +4.01%  '

$ws.Range('D10').Value = '1.815.63'
$ws.Range('E10').Value = '  +5.23%  '

$ws.Range('D11').Value = '''16.82'
$ws.Range('E11').Value = '  +5.39%  '

$ws.Range('D12').Value = '''0.07108'
$ws.Range('E12').Value = '  +3.25%  '

$ws.Range('D13').Value = '''0.6489'
$ws.Range('E13').Value = '  +6.70%  '

$ws.Range('D14').Value = '''84.54'
$ws.Range('E14').Value = '  +9.89%  '

$ws.Range('D15').Value = '''4.711'
$ws.Range('E15').Value = '  +5.30%  '

$ws.Range('D16').Value = '28.823.23'
$ws.Range('E16').Value = '  +7.76%  '

$ws.Range('D17').Value = '''0.9995'
$ws.Range('E17').Value = '  +0.15%  '

$ws.Range('D18').Value = '''0.000007397'
$ws.Range('E18').Value = '  +3.38%  '

$ws.Range('D19').Value = '''0.9991'
$ws.Range('E19').Value = '  +0.20%  '

$ws.Range('E20').Value = '  +7.24%  '

$ws.Range('D21').Value = '2.048.48'
$ws.Range('E21').Value = '  +5.06%  '

$ws.Range('E22').Value = '  +3.90%  '

$ws.Range('D23').Value = '''8.941'
$ws.Range('E23').Value = '  +4.18%  '

$ws.Range('E24').Value = '  +5.34%  '

$ws.Range('D25').Value = '''143.79'
$ws.Range('E25').Value = '  +4.49%  '

$ws.Range('D26').Value = '''133.26'
$ws.Range('E26').Value = '  +25.79%  '

$ws.Range('E27').Value = '  +9.63%  '

$ws.Range('D28').Value = '''1.895'
$ws.Range('E28').Value = '  +6.93%  '

$ws.Range('E29').Value = '  +1.42%  '

$ws.Range('D30').Value = '''4.165'
$ws.Range('E30').Value = '  +4.41%  '

$ws.Range('D31').Value = '''0.08358'
$ws.Range('E31').Value = '  +5.34%  '

$ws.Range('D32').Value = '''3.837'
$ws.Range('E32').Value = '  +3.93%  '

$ws.Range('D33').Value = '''0.04983'
$ws.Range('E33').Value = '  +11.39%  '

$ws.Range('D35').Value = '''0.6779'
$ws.Range('E35').Value = '  +9.20%  '

$ws.Range('D36').Value = '''2.696'
$ws.Range('E36').Value = '  +3.84%  '

$ws.Range('D37').Value = '''2.273'
$ws.Range('E37').Value = '  +11.91%  '

$ws.Range('D38').Value = '''2.755'
$ws.Range('E38').Value = '  +12.97%  '

$ws.Range('D39').Value = '''0.9574'
$ws.Range('E39').Value = '  +3.84%  '

$ws.Range('D40').Value = '''6.051'
$ws.Range('E40').Value = '  +7.10%  '

$ws.Range('D41').Value = '''0.01594'
$ws.Range('E41').Value = '  +6.63%  '

$ws.Range('D42').Value = '''0.9995'
$ws.Range('E42').Value = '  +0.25%  '

$ws.Range('D43').Value = '''0.4098'
$ws.Range('E43').Value = '  +6.60%  '

$ws.Range('D44').Value = '''100.12'
$ws.Range('E44').Value = '  +0.33%  '

$ws.Range('D45').Value = '''7.224'
$ws.Range('E45').Value = '  +5.41%  '

$ws.Range('D46').Value = '''0.1225'
$ws.Range('E46').Value = '  +5.92%  '

$ws.Range('D47').Value = '''0.05532'
$ws.Range('E47').Value = '  +2.76%  '

$ws.Range('D48').Value = '''8.115'
$ws.Range('E48').Value = '  +3.05%  '

$ws.Range('D49').Value = '''31.65'

$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').Value = '''0.3632'
$ws.Range('E50').Value = '  +7.98%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''1.308'
$ws.Range('E51').Value = '  +6.28%  '
